$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look numeric,
# so Excel does not auto-convert them from text to a number (matches source data as text).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "66.696.48"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "3.236.00"
$ws.Range("E3").Value = "  +1.78%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "606.33"
$ws.Range("E5").Value = "  +1.93%  "
$ws.Range("D6").Value = "158.58"
$ws.Range("E6").Value = "  +3.42%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.234.81"
$ws.Range("E8").Value = "  +1.80%  "
$ws.Range("D9").Value = "0.551"
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("E10").Value = "  +2.07%  "
$ws.Range("D11").Value = "5.71"
$ws.Range("E11").Value = "  -5.24%  "
$ws.Range("D12").Value = "0.508"
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("D13").Value = "0.0000276"
$ws.Range("E13").Value = "  +4.14%  "
$ws.Range("D14").Value = "39.10"
$ws.Range("E14").Value = "  +0.74%  "
$ws.Range("D15").Value = "3.764.50"
$ws.Range("E15").Value = "  +1.68%  "
$ws.Range("D16").Value = "66.737.88"
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").Value = "3.232.30"
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("E19").Value = "  +1.18%  "
$ws.Range("D20").Value = "511.98"
$ws.Range("E20").Value = "  +0.62%  "
$ws.Range("D21").Value = "15.24"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D22").Value = "0.736"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("D23").Value = "8.07"
$ws.Range("E23").Value = "  +0.91%  "
$ws.Range("D24").Value = "14.66"
$ws.Range("E24").Value = "  -2.18%  "
$ws.Range("D25").Value = "85.10"
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("D28").Value = "9.20"
$ws.Range("E28").Value = "  +0.65%  "
$ws.Range("D29").Value = "2.40"
$ws.Range("E29").Value = "  +5.13%  "
$ws.Range("D30").Value = "2.97"
$ws.Range("E30").Value = "  +3.05%  "
$ws.Range("D31").Value = "7.04"
$ws.Range("E31").Value = "  +0.88%  "
$ws.Range("D32").Value = "28.29"
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").Value = "1.19"
$ws.Range("E34").Value = "  -2.21%  "
$ws.Range("D35").Value = "0.107"
$ws.Range("E35").Value = "  +19.29%  "
$ws.Range("D36").Value = "6.52"
$ws.Range("E36").Value = "  +0.72%  "
$ws.Range("D37").Value = "510.90"
$ws.Range("E37").Value = "  +4.90%  "
$ws.Range("D38").Value = "55.72"
$ws.Range("E38").Value = "  +1.83%  "
$ws.Range("D39").Value = "0.0₃0780"
$ws.Range("E39").Value = "  +19.71%  "
$ws.Range("E40").Value = "  +0.80%  "
$ws.Range("E41").Value = "  +9.79%  "
$ws.Range("E42").Value = "  +6.18%  "
$ws.Range("E43").Value = "  -1.22%  "
$ws.Range("E44").Value = "  +0.66%  "
$ws.Range("D45").Value = "2.47"
$ws.Range("E45").Value = "  +2.73%  "
$ws.Range("D46").Value = "2.879.01"
$ws.Range("E46").Value = "  -0.51%  "
$ws.Range("D47").Value = "28.56"
$ws.Range("E47").Value = "  +0.95%  "
$ws.Range("D48").Value = "2.43"
$ws.Range("E48").Value = "  +5.56%  "
$ws.Range("D50").Value = "0.117"
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("D51").Value = "122.61"
$ws.Range("E51").Value = "  +0.69%  "
